$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "bom" defined name (was used by the legacy query table range)
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
  $wb.Names.Item($i).Delete()
}

# Convert the old query-table-backed range into a normal Excel Table (ListObject)
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:F29"), $null, 1)
$lo.Name = "Table1"

# Give the "Order link" column the same Hyperlink cell style used by its data
$lo.ListColumns.Item(6).DataBodyRange.Style = "Hyperlink"

# Add the board-cost summary rows below the table
$ws.Range("A32").Value = "Board cost:"
$ws.Range("A33").Formula = "=SUMPRODUCT(C2:C32,E2:E32)"

# Widen column C slightly (no longer auto-fit) to fit the new content
$ws.Columns.Item(3).ColumnWidth = 8.8

# Restore the view: scroll down a bit and select C24 (matches the saved view state)
[void]$ws.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 6
